# Halicarnassus_Profits: scheduled market-data refresh
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H-N)
# for specific leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

function Set-LeveRow {
    param(
        [string]$SheetName,
        [int]$Row,
        [int]$ExpectedItemId,
        [hashtable]$Values
    )
    $ws = $wb.Worksheets.Item($SheetName)
    $actualItemId = $ws.Range("G$Row").Value
    if ($actualItemId -ne $ExpectedItemId) {
        Write-Output "WARNING: $SheetName!G$Row = $actualItemId, expected $ExpectedItemId"
    }
    foreach ($col in $Values.Keys) {
        $val = $Values[$col]
        if ($null -eq $val) {
            $ws.Range("$col$Row").ClearContents()
        } else {
            $ws.Range("$col$Row").Value = $val
        }
    }
}

Set-LeveRow -SheetName "ALC" -Row 4 -ExpectedItemId 5470 -Values @{ "H" = 367.4; "I" = 367.4; "K" = 367.4; "M" = -253.4 }
Set-LeveRow -SheetName "ALC" -Row 5 -ExpectedItemId 5503 -Values @{ "H" = 100.5; "I" = 106.6; "J" = 70; "K" = 106.6; "L" = 70; "M" = 8.400000000000006; "N" = -300 }
Set-LeveRow -SheetName "ALC" -Row 40 -ExpectedItemId 5505 -Values @{ "H" = 6677.136; "I" = 5428.7144; "J" = 8861.875; "K" = 5428.7144; "L" = 8861.875; "M" = -5253.7144; "N" = -9211.875 }
Set-LeveRow -SheetName "ALC" -Row 64 -ExpectedItemId 5506 -Values @{ "H" = 10000; "I" = 0; "K" = 0; "M" = $null }
Set-LeveRow -SheetName "ALC" -Row 67 -ExpectedItemId 5506 -Values @{ "H" = 10000; "I" = 0; "K" = 0; "M" = $null }
Set-LeveRow -SheetName "ALC" -Row 70 -ExpectedItemId 12604 -Values @{ "H" = 4657; "I" = 2118; "K" = 6354; "M" = -6084 }
Set-LeveRow -SheetName "ALC" -Row 73 -ExpectedItemId 12604 -Values @{ "H" = 4657; "I" = 2118; "K" = 6354; "M" = -5418 }
Set-LeveRow -SheetName "ALC" -Row 74 -ExpectedItemId 5507 -Values @{ "H" = 3250; "I" = 3666.6667; "K" = 3666.6667; "M" = -2730.6667 }
Set-LeveRow -SheetName "ALC" -Row 77 -ExpectedItemId 5507 -Values @{ "H" = 3250; "I" = 3666.6667; "K" = 18333.3335; "M" = -13653.3335 }
Set-LeveRow -SheetName "ALC" -Row 106 -ExpectedItemId 19903 -Values @{ "H" = 4122.7144; "I" = 4254.3335; "J" = 3333; "K" = 4254.3335; "L" = 3333; "M" = -3623.3335; "N" = -4595 }
Set-LeveRow -SheetName "ALC" -Row 125 -ExpectedItemId 36228 -Values @{ "H" = 4409.2856; "I" = 4245; "J" = 4628.3335; "K" = 38205; "L" = 41655.0015; "M" = -35745; "N" = -46575.0015 }
Set-LeveRow -SheetName "ARM" -Row 61 -ExpectedItemId 43999 -Values @{ "H" = 3460.0344; "I" = 2573.68; "K" = 2573.68; "M" = -2361.68 }
Set-LeveRow -SheetName "ARM" -Row 110 -ExpectedItemId 27708 -Values @{ "H" = 1387.8; "I" = 1170.625; "K" = 1170.625; "M" = 874.375 }
Set-LeveRow -SheetName "ARM" -Row 122 -ExpectedItemId 36168 -Values @{ "H" = 474.33334; "I" = 490.27274; "K" = 1470.81822; "M" = 979.1817799999999 }
Set-LeveRow -SheetName "ARM" -Row 132 -ExpectedItemId 43997 -Values @{ "H" = 7187.25; "I" = 7187.25; "K" = 21561.75; "M" = -19031.75 }
Set-LeveRow -SheetName "ARM" -Row 136 -ExpectedItemId 43999 -Values @{ "H" = 3460.0344; "I" = 2573.68; "K" = 7721.039999999999; "M" = -5171.039999999999 }
Set-LeveRow -SheetName "BSM" -Row 20 -ExpectedItemId 14149 -Values @{ "H" = 2471; "I" = 1796.1428; "J" = 3061.5; "K" = 1796.1428; "L" = 3061.5; "M" = -1549.1428; "N" = -3555.5 }
Set-LeveRow -SheetName "BSM" -Row 99 -ExpectedItemId 19943 -Values @{ "H" = 1702.4; "I" = 1410.75; "K" = 1410.75; "M" = 87.25 }
Set-LeveRow -SheetName "BSM" -Row 107 -ExpectedItemId 27706 -Values @{ "H" = 2974.9614; "I" = 2058.652; "K" = 2058.652; "M" = -138.652 }
Set-LeveRow -SheetName "CRP" -Row 7 -ExpectedItemId 5361 -Values @{ "H" = 2754.359; "I" = 3924.8845; "J" = 413.30768; "K" = 3924.8845; "L" = 413.30768; "M" = -3811.8845; "N" = -639.30768 }
Set-LeveRow -SheetName "CRP" -Row 16 -ExpectedItemId 27691 -Values @{ "H" = 1961.8; "I" = 1770; "J" = 2249.5; "K" = 1770; "L" = 2249.5; "M" = -1483; "N" = -2823.5 }
Set-LeveRow -SheetName "CRP" -Row 22 -ExpectedItemId 5367 -Values @{ "H" = 1466.5; "I" = 1332.7142; "K" = 1332.7142; "M" = -982.7141999999999 }
Set-LeveRow -SheetName "CRP" -Row 62 -ExpectedItemId 12580 -Values @{ "H" = 3416.6667; "I" = 3416.6667; "K" = 3416.6667; "M" = -2792.6667 }
Set-LeveRow -SheetName "CRP" -Row 65 -ExpectedItemId 12580 -Values @{ "H" = 3416.6667; "I" = 3416.6667; "K" = 17083.3335; "M" = -13963.3335 }
Set-LeveRow -SheetName "CRP" -Row 74 -ExpectedItemId 10636 -Values @{ "H" = 40156.5; "I" = 0; "J" = 40156.5; "K" = 0; "L" = 40156.5; "M" = $null; "N" = -41904.5 }
Set-LeveRow -SheetName "CRP" -Row 77 -ExpectedItemId 10636 -Values @{ "H" = 40156.5; "I" = 0; "J" = 40156.5; "K" = 0; "L" = 120469.5; "M" = $null; "N" = -129205.5 }
Set-LeveRow -SheetName "CRP" -Row 113 -ExpectedItemId 27691 -Values @{ "H" = 1961.8; "I" = 1770; "J" = 2249.5; "K" = 1770; "L" = 2249.5; "M" = 400; "N" = -6589.5 }
Set-LeveRow -SheetName "CUL" -Row 12 -ExpectedItemId 4854 -Values @{ "H" = 104.94118; "I" = 18.75; "J" = 131.46153; "K" = 56.25; "L" = 394.38459; "M" = 116.75; "N" = -740.38459 }
Set-LeveRow -SheetName "CUL" -Row 104 -ExpectedItemId 19807 -Values @{ "H" = 9458.076999999999; "I" = 7500; "K" = 22500; "M" = -19879 }
Set-LeveRow -SheetName "GSM" -Row 102 -ExpectedItemId 36169 -Values @{ "H" = 2674.3333; "I" = 2537.5293; "K" = 2537.5293; "M" = -915.5293000000001 }
Set-LeveRow -SheetName "GSM" -Row 107 -ExpectedItemId 27802 -Values @{ "H" = 1694.4166; "I" = 2159; "J" = 1229.8334; "K" = 2159; "L" = 1229.8334; "M" = -239; "N" = -5069.8334 }
Set-LeveRow -SheetName "GSM" -Row 113 -ExpectedItemId 27710 -Values @{ "H" = 5506; "I" = 2809.8; "J" = 9999.666999999999; "K" = 2809.8; "L" = 9999.666999999999; "M" = -639.8000000000002; "N" = -14339.667 }
Set-LeveRow -SheetName "LTW" -Row 22 -ExpectedItemId 5277 -Values @{ "H" = 688.25; "I" = 449.5; "J" = 927; "K" = 449.5; "L" = 927; "M" = -154.5; "N" = -1517 }
Set-LeveRow -SheetName "LTW" -Row 27 -ExpectedItemId 5277 -Values @{ "H" = 688.25; "I" = 449.5; "J" = 927; "K" = 449.5; "L" = 927; "M" = -342.5; "N" = -1141 }
Set-LeveRow -SheetName "LTW" -Row 40 -ExpectedItemId 36248 -Values @{ "H" = 3209.0908; "J" = 3005; "L" = 3005; "N" = -3277 }
Set-LeveRow -SheetName "LTW" -Row 61 -ExpectedItemId 27740 -Values @{ "H" = 4448.923; "I" = 3274.6667; "J" = 5455.4287; "K" = 3274.6667; "L" = 5455.4287; "M" = -3072.6667; "N" = -5859.4287 }
Set-LeveRow -SheetName "LTW" -Row 100 -ExpectedItemId 19995 -Values @{ "H" = 5202.7334; "I" = 1880.25; "K" = 1880.25; "M" = -1339.25 }
Set-LeveRow -SheetName "LTW" -Row 113 -ExpectedItemId 27740 -Values @{ "H" = 4448.923; "I" = 3274.6667; "J" = 5455.4287; "K" = 3274.6667; "L" = 5455.4287; "M" = -1104.6667; "N" = -9795.4287 }
Set-LeveRow -SheetName "LTW" -Row 116 -ExpectedItemId 26133 -Values @{ "H" = 161666.67; "J" = 161666.67; "L" = 161666.67; "N" = -170844.67 }
Set-LeveRow -SheetName "WVR" -Row 107 -ExpectedItemId 27746 -Values @{ "H" = 2321.5; "I" = 2786.875; "K" = 8360.625; "M" = -6440.625 }
Set-LeveRow -SheetName "WVR" -Row 113 -ExpectedItemId 27752 -Values @{ "H" = 883.4; "I" = 805.6667; "J" = 1000; "K" = 2417.0001; "L" = 3000; "M" = -247.0001000000002; "N" = -7340 }
Set-LeveRow -SheetName "WVR" -Row 122 -ExpectedItemId 36208 -Values @{ "H" = 5950; "I" = 3916.6667; "K" = 11750.0001; "M" = -9300.000100000001 }
Set-LeveRow -SheetName "WVR" -Row 126 -ExpectedItemId 36210 -Values @{ "H" = 4541.4814; "I" = 2895.2942; "J" = 7340; "K" = 8685.882599999999; "L" = 22020; "M" = -6215.882599999999; "N" = -26960 }

Write-Output "Applied scheduled market-data refresh to 43 leve rows across 8 sheets."
